$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New applicant rows to append (rows 152-155), mirroring the existing
# table layout: F.I.Sh | Yo'nalish | Ta'lim tili | Ta'lim shakli |
# Passport | JSHIR | Viloyat | Tuman | Telegram raqami | Telefon raqami | Sana
$rows = @(
    @{ Row = 152; Values = @(
        "Fayzullayev Fatxullo",
        "Yurisprudensiya",
        "O'zbek tili",
        "Kunduzgi",
        "AE1624128",
        "52008075500011",
        "Toshkent shahri",
        "Yashnaobod tumani",
        "998885223533",
        "+998338390777",
        "2025-07-02"
    ) },
    @{ Row = 153; Values = @(
        "Joraqulov Bahodir",
        "Yurisprudensiya",
        "O'zbek tili",
        "Kunduzgi",
        "AD6150839",
        "50608075780064",
        "Navoiy viloyati",
        "Navbahor tumani",
        "998774774313",
        "+998774774313",
        "2025-07-02"
    ) },
    @{ Row = 154; Values = @(
        "Berdiyev Jo'rabek",
        "Yurisprudensiya",
        "O'zbek tili",
        "Kunduzgi",
        "AB6139569",
        "50107005800024",
        "Jizzax viloyati",
        "Sharof Rashidov tumani",
        "998992428774",
        "+998992428774",
        "2025-07-02"
    ) },
    @{ Row = 155; Values = @(
        "Shorustamova Sabina Shoakbar qizi",
        "Yurisprudensiya",
        "Rus tili",
        "Kunduzgi",
        "AD5424978",
        "62209076510013",
        "Toshkent shahri",
        "Bektemir tumani",
        "998999290005",
        "+998999290005",
        "2025-07-02"
    ) }
)

# Use an existing plain data row as the style template so the appended
# rows stay unstyled (no explicit cell style index), same as the other
# rows in the table.
$styleTemplate = $ws.Range("A3:K3")

foreach ($entry in $rows) {
    $r = $entry.Row
    $rowRange = $ws.Range("A" + $r + ":K" + $r)

    # Force text formatting first so numeric-looking values (passport
    # series/numbers, JSHIR, phone numbers, the ISO date string) are
    # stored as text instead of being auto-coerced to numbers/dates.
    $rowRange.NumberFormat = "@"

    for ($i = 0; $i -lt $entry.Values.Length; $i++) {
        $cell = $ws.Cells.Item($r, $i + 1)
        $cell.Value = $entry.Values[$i]
    }

    # Re-apply the plain/unstyled look of the surrounding data rows.
    $rowRange.Style = $styleTemplate.Style
}
